$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata" ---
$meta = $wb.Worksheets.Item("Metadata")

# Update Version value (row 3)
$meta.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"

# Update Date value (row 8)
$meta.Range("B8").Value = "2025-10-29T22:15:57+01:00"

# Insert a new row 11 for "Jurisdiction", shifting existing rows (old 11-19) down to 12-20
$meta.Rows.Item(11).Insert()

# Copy formatting from the row now below (row 12, the old row 11 "Description") so the
# new row matches the sheet's standard data-row style
$meta.Range("A12:B12").Copy()
$meta.Range("A11:B11").PasteSpecial(-4122)

# Fill in the new Jurisdiction property row (no value)
$meta.Range("A11").Value = "Jurisdiction"
$meta.Range("B11").Value = ""

# --- Sheet "Elements" ---
$elements = $wb.Worksheets.Item("Elements")

# LegalAuthenticator.typeId (row 5) gains an invariant constraint on its II type
$elements.Range("AJ5").Value = "II-1:An II instance must have either a root or an nullFlavor. {root.exists() or nullFlavor.exists()}
"
